$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland I Liga")

# Rows 5-8 and 86-87 were reordered/rearranged (same matches, different row
# positions) as part of the source data refresh. Apply the final values for
# every affected cell directly.

# Row 5: Chojniczanka Chojnice vs GKS Katowice
$ws.Range("B5").Value = 5448050
$ws.Range("E5").Value = "Chojniczanka Chojnice"
$ws.Range("F5").Value = "GKS Katowice"
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = "D"
$ws.Range("L5").Value = 2.75
$ws.Range("M5").Value = 3.25
$ws.Range("N5").Value = 2.375
$ws.Range("O5").Value = 2.4
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 1.8
$ws.Range("T5").Value = 2.05
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 1.85
$ws.Range("X5").Value = -1
$ws.Range("Y5").Value = 2.25
$ws.Range("Z5").Value = -1
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 1
$ws.Range("AD5").Value = -1

# Row 6: Zaglebie Sosnowiec vs Skra Czestochowa
$ws.Range("B6").Value = 5448048
$ws.Range("E6").Value = "Zaglebie Sosnowiec"
$ws.Range("F6").Value = "Skra Czestochowa"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "D"
$ws.Range("L6").Value = 2.1
$ws.Range("M6").Value = 3.2
$ws.Range("N6").Value = 3.3
$ws.Range("O6").Value = 2.1
$ws.Range("P6").Value = 3.2
$ws.Range("Q6").Value = 3.1
$ws.Range("R6").Value = -0.25
$ws.Range("S6").Value = 1.875
$ws.Range("T6").Value = 1.925
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.85
$ws.Range("W6").Value = 1.95
$ws.Range("X6").Value = -1
$ws.Range("Y6").Value = 2.2
$ws.Range("Z6").Value = -1
$ws.Range("AA6").Value = -0.5
$ws.Range("AB6").Value = 0.4625
$ws.Range("AC6").Value = -0.5
$ws.Range("AD6").Value = 0.475

# Row 7: Stal Rzeszow vs Sandecja Nowy Sacz
$ws.Range("B7").Value = 5451610
$ws.Range("E7").Value = "Stal Rzeszow"
$ws.Range("F7").Value = "Sandecja Nowy Sacz"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = "H"
$ws.Range("L7").Value = 1.444
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 6.5
$ws.Range("O7").Value = 1.333
$ws.Range("P7").Value = 4.333
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = -1.5
$ws.Range("S7").Value = 1.95
$ws.Range("T7").Value = 1.85
$ws.Range("U7").Value = 2.75
$ws.Range("V7").Value = 1.875
$ws.Range("W7").Value = 1.925
$ws.Range("X7").Value = 0.333
$ws.Range("Y7").Value = -1
$ws.Range("Z7").Value = -1
$ws.Range("AA7").Value = -1
$ws.Range("AB7").Value = 0.8500000000000001
$ws.Range("AC7").Value = 0.4375
$ws.Range("AD7").Value = -0.5

# Row 8: LKS Lodz vs Odra Opole
$ws.Range("B8").Value = 5448049
$ws.Range("E8").Value = "LKS Lodz"
$ws.Range("F8").Value = "Odra Opole"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = "H"
$ws.Range("L8").Value = 1.571
$ws.Range("M8").Value = 3.75
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 1.444
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = -1
$ws.Range("S8").Value = 1.775
$ws.Range("T8").Value = 2.025
$ws.Range("U8").Value = 2.75
$ws.Range("V8").Value = 1.9
$ws.Range("W8").Value = 1.9
$ws.Range("X8").Value = 0.444
$ws.Range("Y8").Value = -1
$ws.Range("Z8").Value = -1
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = -1
$ws.Range("AD8").Value = 0.8999999999999999

# Row 86: Miedz Legnica vs Odra Opole
$ws.Range("B86").Value = 6803740
$ws.Range("E86").Value = "Miedz Legnica"
$ws.Range("F86").Value = "Odra Opole"
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 2
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1
$ws.Range("K86").Value = "A"
$ws.Range("L86").Value = 1.85
$ws.Range("M86").Value = 3.5
$ws.Range("N86").Value = 3.75
$ws.Range("O86").Value = 1.909
$ws.Range("P86").Value = 3.5
$ws.Range("Q86").Value = 3.5
$ws.Range("R86").Value = -0.5
$ws.Range("S86").Value = 1.975
$ws.Range("T86").Value = 1.825
$ws.Range("U86").Value = 2.25
$ws.Range("V86").Value = 1.9
$ws.Range("W86").Value = 1.9
$ws.Range("X86").Value = -1
$ws.Range("Y86").Value = -1
$ws.Range("Z86").Value = 2.5
$ws.Range("AA86").Value = -1
$ws.Range("AB86").Value = 0.825
$ws.Range("AC86").Value = 0.8999999999999999
$ws.Range("AD86").Value = -1

# Row 87: Podbeskidzie Bielsko Biala vs Gornik Leczna
$ws.Range("B87").Value = 6803738
$ws.Range("E87").Value = "Podbeskidzie Bielsko Biala"
$ws.Range("F87").Value = "Gornik Leczna"
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 1
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = "D"
$ws.Range("L87").Value = 1.85
$ws.Range("M87").Value = 3.5
$ws.Range("N87").Value = 3.8
$ws.Range("O87").Value = 1.666
$ws.Range("P87").Value = 3.8
$ws.Range("Q87").Value = 4.75
$ws.Range("R87").Value = -0.75
$ws.Range("S87").Value = 1.825
$ws.Range("T87").Value = 1.975
$ws.Range("U87").Value = 2.5
$ws.Range("V87").Value = 1.825
$ws.Range("W87").Value = 1.975
$ws.Range("X87").Value = -1
$ws.Range("Y87").Value = 2.8
$ws.Range("Z87").Value = -1
$ws.Range("AA87").Value = -1
$ws.Range("AB87").Value = 0.9750000000000001
$ws.Range("AC87").Value = -1
$ws.Range("AD87").Value = 0.9750000000000001

